# feat: Rules checks and JSON simplification
#
# 1) Add a new "listed_simple" sheet after "listed" - a simplified version
#    of the "listed" table (id + a plain integer column instead of the
#    nominal/ordinal/date columns).
# 2) Normalise the cell formatting on "listed" (and apply the same
#    formatting to the new sheet) so both share one explicit style.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Normalise formatting on the existing "listed" sheet.
# ---------------------------------------------------------------------
$listed = $wb.Worksheets.Item("listed")
$listedRange = $listed.Range("A1:C12")
$listedRange.Locked = $false

# Move the selection off the old active cell (C1) to A12, matching the
# post-edit workbook state.
$listed.Range("A12").Select()

# ---------------------------------------------------------------------
# 2) Add the new "listed_simple" sheet right after "listed".
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$simple = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$simple.Name = "listed_simple"

# Header row.
$simple.Range("A1").Value = "id"
$simple.Range("B1").Value = "integer"

# Data rows - same "id" grouping as "listed", paired with a plain integer.
$ids = @("id01", "id01", "id02", "id03", "id03", "id03", "id04", "id05", "id05", "id05", "id05")
$ints = @(234, 22, 54, 34, 1, 54, 76, 23, 45, 23, 65)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $simple.Cells.Item($row, 1).Value = $ids[$i]
    $simple.Cells.Item($row, 2).Value = $ints[$i]
}

# Same explicit-but-invisible formatting as "listed", so both sheets share
# the same style.
$simple.Range("A1:B12").Locked = $false

# The new sheet becomes the active tab, selection on the cell just below
# the data.
$simple.Range("B13").Select()
